$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 464
$ws.Range("F4").Value = 492
$ws.Range("F5").Value = 2282
$ws.Range("F7").Value = 61
$ws.Range("F9").Value = 1674
$ws.Range("F10").Value = 1674
$ws.Range("F11").Value = 1382
$ws.Range("F12").Value = 69
$ws.Range("F13").Value = 1432
$ws.Range("F15").Value = 22
$ws.Range("F16").Value = 645
$ws.Range("F17").Value = 174
$ws.Range("F18").Value = 122
$ws.Range("F19").Value = 7389
$ws.Range("F20").Value = 8244
$ws.Range("F23").Value = 214
$ws.Range("F25").Value = 492
$ws.Range("F27").Value = 220
$ws.Range("F28").Value = 267
$ws.Range("F29").Value = 255
$ws.Range("F30").Value = 19
$ws.Range("F31").Value = 2
$ws.Range("F33").Value = 354
$ws.Range("F34").Value = 1476
$ws.Range("F35").Value = 252
$ws.Range("F36").Value = 235
$ws.Range("F38").Value = 297
$ws.Range("F39").Value = 27
$ws.Range("F40").Value = 761
$ws.Range("F41").Value = 25
$ws.Range("F42").Value = 1369
$ws.Range("F43").Value = 358
$ws.Range("F44").Value = 257
$ws.Range("F45").Value = 209
$ws.Range("F46").Value = 92
$ws.Range("F47").Value = 192
$ws.Range("F48").Value = 176
$ws.Range("F49").Value = 15
$ws.Range("F50").Value = 17

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 35
$ws.Range("F13").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 188
$ws.Range("F3").Value = 2638
$ws.Range("F4").Value = 291
$ws.Range("F5").Value = 144

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 188
$ws.Range("F4").Value = 464
$ws.Range("F6").Value = 144
$ws.Range("F7").Value = 35
$ws.Range("F8").Value = 492
$ws.Range("F9").Value = 2282
$ws.Range("F10").Value = 61
$ws.Range("F12").Value = 1674
$ws.Range("F13").Value = 1674
$ws.Range("F14").Value = 69
$ws.Range("F15").Value = 22
$ws.Range("F16").Value = 645
$ws.Range("F18").Value = 174
$ws.Range("F20").Value = 122
$ws.Range("F21").Value = 7389
$ws.Range("F22").Value = 8244
$ws.Range("F24").Value = 214
$ws.Range("F26").Value = 267
$ws.Range("F27").Value = 19
$ws.Range("F29").Value = 252
$ws.Range("F30").Value = 235
$ws.Range("F33").Value = 297
$ws.Range("F34").Value = 27
$ws.Range("F37").Value = 761
$ws.Range("F39").Value = 25
$ws.Range("F40").Value = 11
$ws.Range("F42").Value = 1369
$ws.Range("F43").Value = 358
$ws.Range("F44").Value = 257
$ws.Range("F45").Value = 209
$ws.Range("F46").Value = 92
$ws.Range("F47").Value = 192
$ws.Range("F50").Value = 17
